$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Pressing opposite button switches indicator." - merge the 3 runs that
#    make up this sentence into a single run by doing a same-text find and
#    replace across the whole sentence (Word rewrites the range as one run).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Pressing opposite button switches indicator.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Pressing opposite button switches indicator.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "UART Logging of: " - merge the bold "UART Logging of:" run and the
#    trailing bold space run into a single bold run with a preserved space.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "UART Logging of: ", $true, $false, $false, $false, $false, $true, 1,
    $false, "UART Logging of: ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "UART logs saved via Tera Term are included in project submission."
#    - drop the bullet numbering, indent the paragraph instead, and wrap the
#      sentence in literal parentheses held in their own runs.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "UART logs saved via Tera Term are included in project submission.*") {
        $frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr>' +
                '<w:r><w:t>(</w:t></w:r>' +
                '<w:r><w:t>UART logs saved via Tera Term are included in project submission.</w:t></w:r>' +
                '<w:r><w:t>)</w:t></w:r>' +
                '</w:p>'
        $p.Range.InsertXML($frag)
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Turn the plain-text GitHub repository URL into a real hyperlink.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "https://github.com/vivekanandaramanu/esp32-indicator-system*" -and
        $p.Range.Text -notlike "*blob*") {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $d.Hyperlinks.Add($r, "https://github.com/vivekanandaramanu/esp32-indicator-system", "", "",
            "https://github.com/vivekanandaramanu/esp32-indicator-system") | Out-Null
        # Materialise the "Hyperlink" character style definition into
        # styles.xml (Hyperlinks.Add alone only references it by name).
        $r.Style = "Hyperlink"
        break
    }
}

# ---------------------------------------------------------------------------
# 5) Turn the plain-text Google Drive URL into a real hyperlink.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "https://drive.google.com/file/d/1q4hKbzImOk4m6uuVi8Ycv6JALjvGyG4a/view?usp=sharing*") {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null
        $d.Hyperlinks.Add($r, "https://drive.google.com/file/d/1q4hKbzImOk4m6uuVi8Ycv6JALjvGyG4a/view?usp=sharing", "", "",
            "https://drive.google.com/file/d/1q4hKbzImOk4m6uuVi8Ycv6JALjvGyG4a/view?usp=sharing") | Out-Null
        $r.Style = "Hyperlink"
        break
    }
}

# ---------------------------------------------------------------------------
# 6) After "Log captured using HW-417-V1.2 USB to TTL module and Tera Term."
#    insert three new paragraphs: a hyperlink to the log file on GitHub, an
#    "OR" separator, and a hyperlink back to the repository root.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Log captured using HW-417-V1.2 USB to TTL module and Tera Term.*") {
        $r = $p.Range
        $r.Collapse(0) | Out-Null
        $r.InsertParagraphAfter() | Out-Null

        $p1 = $p.Next()
        $p1.Range.InsertAfter("https://github.com/vivekanandaramanu/esp32-indicator-system/blob/main/indicator_log.txt") | Out-Null
        $rng1 = $p1.Range
        $rng1.MoveEnd(1, -1) | Out-Null
        $d.Hyperlinks.Add($rng1, "https://github.com/vivekanandaramanu/esp32-indicator-system/blob/main/indicator_log.txt", "", "",
            "https://github.com/vivekanandaramanu/esp32-indicator-system/blob/main/indicator_log.txt") | Out-Null
        $rng1.Style = "Hyperlink"

        $p1.Range.InsertParagraphAfter() | Out-Null
        $p2 = $p1.Next()
        $p2.Range.InsertAfter("OR") | Out-Null

        $p2.Range.InsertParagraphAfter() | Out-Null
        $p3 = $p2.Next()
        $p3.Range.InsertAfter("https://github.com/vivekanandaramanu/esp32-indicator-system") | Out-Null
        $rng3 = $p3.Range
        $rng3.MoveEnd(1, -1) | Out-Null
        $d.Hyperlinks.Add($rng3, "https://github.com/vivekanandaramanu/esp32-indicator-system", "", "",
            "https://github.com/vivekanandaramanu/esp32-indicator-system") | Out-Null
        $rng3.Style = "Hyperlink"

        break
    }
}

# ---------------------------------------------------------------------------
# 7) Style metadata touch-ups mirroring what Word stamps into styles.xml the
#    first time a Hyperlink (and the related Unresolved Mention) style is
#    used in a session: based on Default Paragraph Font, visible in the
#    gallery, classic hyperlink blue with the hyperlink theme color.
#    NOTE: the Hyperlink character style was already materialised into
#    styles.xml by the Hyperlinks.Add() calls above, so it is now safe to
#    tweak its properties without the edits leaking onto unrelated runs.
# ---------------------------------------------------------------------------
$hs = $d.Styles("Hyperlink")
$hs.Priority = 99
$hs.UnhideWhenUsed = $true
$hs.Font.Color = 16711680

try {
    $um = $d.Styles.Add("Unresolved Mention", 2)
    $um.BaseStyle = $d.Styles("DefaultParagraphFont")
    $um.Priority = 99
    $um.UnhideWhenUsed = $true
} catch {
}
